$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# --- Elements sheet: fix casing of "ExerciceProfessionnel" references ---
$wsElem = $wb.Worksheets.Item("Elements")

# Row 6 is the FonctionQualifiee.exerciceProfessionnel element
$wsElem.Range("A6").Value = "FonctionQualifiee.ExerciceProfessionnel"
$wsElem.Range("B6").Value = "FonctionQualifiee.ExerciceProfessionnel"
$wsElem.Range("L6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("M6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("AF6").Value = "SavoirFaire.ExerciceProfessionnel"

# Column A/B best-fit width grows slightly wider to accommodate the new text
# (target stored width 32.30078125; engine quantizes ColumnWidth to whole
# pixels, so 31.5 is the closest input that lands on the nearest reachable
# stored width)
$wsElem.Columns.Item(1).ColumnWidth = 31.5
$wsElem.Columns.Item(2).ColumnWidth = 31.5
